$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.704275
$ws.Range("H2").Value = 20.112825
$ws.Range("I2").Value = 0.4617710489234531
$ws.Range("J2").Value = 0.4617710489234532
$ws.Range("M2").Value = 25.37147633333333
$ws.Range("N2").Value = 76.114429
$ws.Range("O2").Value = 0.5780881462719274
$ws.Range("P2").Value = 0.5780881462719274
$ws.Range("Q2").Value = 170.0973544946583
$ws.Range("R2").Value = 1530.876190451925
$ws.Range("S2").Value = 0.2669443696742025
$ws.Range("T2").Value = 0.2669443696742025
$ws.Range("G3").Value = 6.704275
$ws.Range("H3").Value = 20.112825
$ws.Range("I3").Value = 0.4617710489234531
$ws.Range("J3").Value = 0.4617710489234532
$ws.Range("O3").Value = 0.2328552951919536
$ws.Range("P3").Value = 0.2328552951919536
$ws.Range("Q3").Value = 68.51562334854168
$ws.Range("R3").Value = 616.640610136875
$ws.Range("S3").Value = 0.1075258339081687
$ws.Range("T3").Value = 0.1075258339081687
$ws.Range("G4").Value = 6.704275
$ws.Range("H4").Value = 20.112825
$ws.Range("I4").Value = 0.4617710489234531
$ws.Range("J4").Value = 0.4617710489234532
$ws.Range("O4").Value = 0.189056558536119
$ws.Range("P4").Value = 0.189056558536119
$ws.Range("Q4").Value = 55.62823016566666
$ws.Range("R4").Value = 500.654071491
$ws.Range("S4").Value = 0.08730084534108187
$ws.Range("T4").Value = 0.08730084534108189
$ws.Range("I5").Value = 0.03922895479591048
$ws.Range("J5").Value = 0.03922895479591048
$ws.Range("M5").Value = 25.37147633333333
$ws.Range("N5").Value = 76.114429
$ws.Range("O5").Value = 0.5780881462719274
$ws.Range("P5").Value = 0.5780881462719274
$ws.Range("Q5").Value = 14.45032434565
$ws.Range("R5").Value = 130.05291911085
$ws.Range("S5").Value = 0.02267779375815313
$ws.Range("T5").Value = 0.02267779375815313
$ws.Range("I6").Value = 0.03922895479591048
$ws.Range("J6").Value = 0.03922895479591048
$ws.Range("O6").Value = 0.2328552951919536
$ws.Range("P6").Value = 0.2328552951919536
$ws.Range("S6").Value = 0.009134669849073539
$ws.Range("T6").Value = 0.009134669849073539
$ws.Range("I7").Value = 0.03922895479591048
$ws.Range("J7").Value = 0.03922895479591048
$ws.Range("O7").Value = 0.189056558536119
$ws.Range("P7").Value = 0.189056558536119
$ws.Range("S7").Value = 0.007416491188683815
$ws.Range("T7").Value = 0.007416491188683816
$ws.Range("I8").Value = 0.4989999962806363
$ws.Range("J8").Value = 0.4989999962806364
$ws.Range("M8").Value = 25.37147633333333
$ws.Range("N8").Value = 76.114429
$ws.Range("O8").Value = 0.5780881462719274
$ws.Range("P8").Value = 0.5780881462719274
$ws.Range("Q8").Value = 183.8109588248586
$ws.Range("R8").Value = 1654.298629423727
$ws.Range("S8").Value = 0.2884659828395718
$ws.Range("T8").Value = 0.2884659828395718
$ws.Range("I9").Value = 0.4989999962806363
$ws.Range("J9").Value = 0.4989999962806364
$ws.Range("O9").Value = 0.2328552951919536
$ws.Range("P9").Value = 0.2328552951919536
$ws.Range("S9").Value = 0.1161947914347113
$ws.Range("T9").Value = 0.1161947914347113
$ws.Range("I10").Value = 0.4989999962806363
$ws.Range("J10").Value = 0.4989999962806364
$ws.Range("O10").Value = 0.189056558536119
$ws.Range("P10").Value = 0.189056558536119
$ws.Range("S10").Value = 0.09433922200635328
$ws.Range("T10").Value = 0.0943392220063533
